# Updated cryptos list snapshot (prices + 1h volume deltas).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range('D2').Value = '26.921.90'

# Row 3: Ethereum
$ws.Range('D3').Value = '1.549.63'
$ws.Range('E3').Value = '  -0.08%  '

# Row 4: TetherUSD
$ws.Range('E4').Value = '  -0.10%  '

# Row 5: BNB
$ws.Range('D5').Value = '''206.70'
$ws.Range('E5').Value = '  +0.28%  '

# Row 6: XRP
$ws.Range('D6').Value = '''0.487'
$ws.Range('E6').Value = '  +0.45%  '

# Row 7: USDC
$ws.Range('E7').Value = '  -0.12%  '

# Row 8: Solana
$ws.Range('D8').Value = '''22.14'
$ws.Range('E8').Value = '  +3.07%  '

# Row 9: Cardano
$ws.Range('E9').Value = '  -0.69%  '

# Row 10: Dogecoin
$ws.Range('E10').Value = '  +0.68%  '

# Row 11: TRON
$ws.Range('D11').Value = '''0.0857'
$ws.Range('E11').Value = '  -0.17%  '

# Row 12: WrappedliquidstakedEther2.0
$ws.Range('D12').Value = '1.769.81'
$ws.Range('E12').Value = '  -0.05%  '

# Row 13: WrappedEther
$ws.Range('D13').Value = '1.549.02'
$ws.Range('E13').Value = '  -0.13%  '

# Row 14: Polkadot
$ws.Range('E14').Value = '  +0.88%  '

# Row 15: Polygon
$ws.Range('E15').Value = '  +1.17%  '

# Row 16: WrappedBTC
$ws.Range('D16').Value = '26.907.80'
$ws.Range('E16').Value = '  -0.13%  '

# Row 17: Litecoin
$ws.Range('D17').Value = '''61.66'
$ws.Range('E17').Value = '  +0.09%  '

# Row 18: BitcoinCash
$ws.Range('D18').Value = '''217.66'
$ws.Range('E18').Value = '  +1.25%  '

# Row 19: ShibaInu
$ws.Range('E19').Value = '  +1.60%  '

# Row 20: Chainlink
$ws.Range('E20').Value = '  +0.56%  '

# Row 22: Uniswap
$ws.Range('E22').Value = '  +0.52%  '

# Row 23: Avalanche
$ws.Range('D23').Value = '''9.19'
$ws.Range('E23').Value = '  -0.03%  '

# Row 24: Toncoin
$ws.Range('E24').Value = '  +0.70%  '

# Row 25: Monero
$ws.Range('D25').Value = '''154.39'
$ws.Range('E25').Value = '  +0.61%  '

# Row 26: Cosmos
$ws.Range('E26').Value = '  -0.53%  '

# Row 27: EthereumClassic
$ws.Range('E27').Value = '  +0.42%  '

# Row 28: Stellar
$ws.Range('E28').Value = '  +0.83%  '

# Row 29: BinanceUSD
$ws.Range('E29').Value = '  -0.09%  '

# Row 30: Hedera
$ws.Range('E30').Value = '  +1.47%  '

# Row 31: PancakeSwap
$ws.Range('E31').Value = '  -0.65%  '

# Row 32: Filecoin
$ws.Range('D32').Value = '''3.23'
$ws.Range('E32').Value = '  -0.17%  '

# Row 33: Maker
$ws.Range('D33').Value = '1.417.84'
$ws.Range('E33').Value = '  +3.34%  '

# Row 34: InternetComputer(DFINITY)
$ws.Range('E34').Value = '  +4.12%  '

# Row 35: LidoDAOToken
$ws.Range('D35').Value = '''1.57'
$ws.Range('E35').Value = '  +2.06%  '

# Row 36: TrustWalletToken
$ws.Range('E36').Value = '  +0.33%  '

# Row 37: HuobiToken
$ws.Range('E37').Value = '  +0.20%  '

# Row 38: VeChain
$ws.Range('E38').Value = '  +0.50%  '

# Row 39: ImmutableX
$ws.Range('D39').Value = '''0.524'
$ws.Range('E39').Value = '  +0.83%  '

# Row 40: ARBITRUM
$ws.Range('D40').Value = '''0.808'
$ws.Range('E40').Value = '  +0.13%  '

# Row 41: FraxShare
$ws.Range('D41').Value = '''5.76'
$ws.Range('E41').Value = '  +5.12%  '

# Row 42: PaxDollar
$ws.Range('E42').Value = '  -0.17%  '

# Row 43: MXToken
$ws.Range('E43').Value = '  +1.05%  '

# Row 44: WEMIXToken
$ws.Range('D44').Value = '''0.993'
$ws.Range('E44').Value = '  +0.76%  '

# Row 45: Aave
$ws.Range('D45').Value = '''64.40'
$ws.Range('E45').Value = '  +1.11%  '

# Row 46: RenderToken
$ws.Range('E46').Value = '  +0.30%  '

# Row 47: RocketPoolETH
$ws.Range('D47').Value = '1.683.99'
$ws.Range('E47').Value = '  -0.01%  '

# Row 48: Quant
$ws.Range('D48').Value = '''87.73'
$ws.Range('E48').Value = '  +1.62%  '

# Row 49: Cronos
$ws.Range('E49').Value = '  +2.17%  '

# Row 50: BabyDogeCoin
$ws.Range('E50').Value = '  +5.45%  '

# Row 51: Algorand
$ws.Range('D51').Value = '''0.0953'
$ws.Range('E51').Value = '  +0.04%  '
